$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Update Price (D) and Volume(1h) (E) columns for rows 2-43 ---
$ws.Range("D2").Value = '29.335.41'
$ws.Range("E2").Value = '  +0.78%  '
$ws.Range("D3").Value = '1.848.41'
$ws.Range("E3").Value = '  +0.74%  '
Set-TextValue $ws.Range("D4") '1.015'
$ws.Range("E4").Value = '  +0.90%  '
Set-TextValue $ws.Range("D5") '244.43'
$ws.Range("E5").Value = '  +0.29%  '
Set-TextValue $ws.Range("D6") '0.6199'
$ws.Range("E6").Value = '  -1.39%  '
Set-TextValue $ws.Range("D7") '1.012'
$ws.Range("E7").Value = '  +0.77%  '
Set-TextValue $ws.Range("D8") '0.07452'
$ws.Range("E8").Value = '  -0.22%  '
Set-TextValue $ws.Range("D9") '0.2954'
$ws.Range("E9").Value = '  +0.84%  '
Set-TextValue $ws.Range("D10") '23.04'
$ws.Range("E10").Value = '  +0.20%  '
Set-TextValue $ws.Range("D11") '0.07732'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = '1.855.63'
$ws.Range("E12").Value = '  +1.44%  '
Set-TextValue $ws.Range("D13") '5.027'
$ws.Range("E13").Value = '  +0.38%  '
Set-TextValue $ws.Range("D14") '0.6738'
$ws.Range("E14").Value = '  +0.95%  '
Set-TextValue $ws.Range("D15") '83.35'
$ws.Range("E15").Value = '  +0.10%  '
Set-TextValue $ws.Range("D16") '0.000009065'
$ws.Range("E16").Value = '  -3.70%  '
Set-TextValue $ws.Range("D17") '5.911'
$ws.Range("E17").Value = '  -2.18%  '
$ws.Range("D18").Value = '29.315.95'
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("D19").Value = '2.070.51'
$ws.Range("E19").Value = '  -0.04%  '
Set-TextValue $ws.Range("D20") '237.74'
$ws.Range("E20").Value = '  +6.17%  '
Set-TextValue $ws.Range("D21") '12.71'
$ws.Range("E21").Value = '  +0.79%  '
Set-TextValue $ws.Range("D22") '1.014'
$ws.Range("E22").Value = '  +0.91%  '
Set-TextValue $ws.Range("D23") '7.196'
$ws.Range("E23").Value = '  +1.23%  '
Set-TextValue $ws.Range("D24") '1.011'
$ws.Range("E24").Value = '  +0.44%  '
Set-TextValue $ws.Range("D25") '160.21'
$ws.Range("E25").Value = '  +0.09%  '
Set-TextValue $ws.Range("D26") '0.1441'
$ws.Range("E26").Value = '  +2.82%  '
Set-TextValue $ws.Range("D27") '8.540'
$ws.Range("E27").Value = '  +0.56%  '
Set-TextValue $ws.Range("D28") '17.96'
$ws.Range("E28").Value = '  +0.16%  '
Set-TextValue $ws.Range("D29") '1.508'
$ws.Range("E29").Value = '  +0.70%  '
Set-TextValue $ws.Range("D30") '4.172'
$ws.Range("E30").Value = '  +0.88%  '
Set-TextValue $ws.Range("D31") '0.05608'
$ws.Range("E31").Value = '  +2.86%  '
Set-TextValue $ws.Range("D32") '4.120'
$ws.Range("E32").Value = '  +1.41%  '
Set-TextValue $ws.Range("D33") '1.218'
$ws.Range("E33").Value = '  +1.42%  '
Set-TextValue $ws.Range("D34") '0.7553'
$ws.Range("E34").Value = '  +0.67%  '
Set-TextValue $ws.Range("D35") '1.856'
$ws.Range("E35").Value = '  +0.36%  '
Set-TextValue $ws.Range("D36") '1.145'
$ws.Range("E36").Value = '  +0.71%  '
Set-TextValue $ws.Range("D37") '2.667'
$ws.Range("E37").Value = '  +2.22%  '
Set-TextValue $ws.Range("D38") '2.843'
$ws.Range("E38").Value = '  +2.98%  '
Set-TextValue $ws.Range("D39") '0.01787'
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("D40").Value = '1.215.20'
$ws.Range("E40").Value = '  -1.14%  '
Set-TextValue $ws.Range("D41") '6.450'
$ws.Range("E41").Value = '  -2.78%  '
Set-TextValue $ws.Range("D42") '0.9101'
$ws.Range("E42").Value = '  +1.71%  '
Set-TextValue $ws.Range("D43") '1.012'
$ws.Range("E43").Value = '  +0.75%  '

# --- Rows 44-51: re-ranked coins (names/links/prices/volumes) ---
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '1.994.51'
$ws.Range("E44").Value = '  +1.16%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D45") '101.16'
$ws.Range("E45").Value = '  -0.42%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D46") '65.47'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D47") '0.00000000123'
$ws.Range("E47").Value = '  +0.59%  '
Set-TextValue $ws.Range("D48") '0.5155'
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D49") '9.228'
$ws.Range("E49").Value = '  +2.76%  '
$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D50") '0.4068'
$ws.Range("E50").Value = '  +0.70%  '
Set-TextValue $ws.Range("D51") '0.05854'
$ws.Range("E51").Value = '  +0.65%  '
